# Se agrega el campo para EGMA
# Adds a new "egma_items" column (column W) as the last header in each of
# the four grade sheets ("2°", "3°", "4°", "5°"), mirroring the manual
# Excel edit described by the commit.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2°", "3°", "4°", "5°")
$selections = @{
    "2°" = "W1"
    "3°" = "W1"
    "4°" = "U6"
    "5°" = "Q8"
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate() | Out-Null
    $ws.Range("W1").Value = "egma_items"
    $ws.Range($selections[$name]).Select() | Out-Null
}

# Leave the workbook with "5°" as the active sheet, matching the final
# state captured in the saved file.
$wb.Worksheets.Item("5°").Activate() | Out-Null
